$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.044.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "'3.765.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'632.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").Value = "'165.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "'3.764.19"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "'6.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  -3.69%  "
$ws.Range("D14").Value = "'34.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "'4.398.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "'3.767.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "'69.021.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "'17.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'7.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "'461.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").Value = "'0.706"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("D25").Value = "'82.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'3.914.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("D32").Value = "'2.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("D34").Value = "'28.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  +18.10%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'3.718.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").Value = "'8.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "'3.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("D41").Value = "'5.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").Value = "'0.965"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D45").Value = "'157.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("D46").Value = "'1.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.53%  "
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").Value = "'47.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'42.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "'0.295"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("E51").Value = "  -0.31%  "
